$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Hasil Analisa" column (F) - its header plus the empty, yellow-highlighted
# cells below it - is no longer needed: select it and clear both contents and
# formatting, same as the user selecting F1:F9 and pressing Delete.
$ws.Range("F1:F9").Select() | Out-Null
$ws.Range("F1:F9").Clear() | Out-Null

# Correct the customer feedback value in C3: "Ada" -> "Ga Ada"
$ws.Range("C3").Value() = "Ga Ada"
